# Refresh the forecast-match statistics (matchonfcst) with the latest
# scheduled run data. Updates the numeric metric columns (B:L) for every
# queue/forecast row on the sheet; labels in column A are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 0
$ws.Cells.Item(2, 2).Value = 3653
$ws.Cells.Item(2, 3).Value = 30
$ws.Cells.Item(2, 4).Value = 650
$ws.Cells.Item(2, 5).Value = 1656
$ws.Cells.Item(2, 6).Value = 145
$ws.Cells.Item(2, 7).Value = 8
$ws.Cells.Item(2, 8).Value = 3746
$ws.Cells.Item(2, 9).Value = 22119
$ws.Cells.Item(2, 10).Value = 17
$ws.Cells.Item(2, 11).Value = 825
$ws.Cells.Item(2, 12).Value = 18065

# Row 3: COV-GESTIONE-COMM-CMN-RES
$ws.Cells.Item(3, 4).Value = 1.322
$ws.Cells.Item(3, 5).Value = 1.327
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(3, 7).Value = 3
$ws.Cells.Item(3, 8).Value = 15
$ws.Cells.Item(3, 9).Value = 8279
$ws.Cells.Item(3, 10).Value = 28

# Row 4: COV-GESTIONE-FO-CT
$ws.Cells.Item(4, 2).Value = 267

# Row 5: COV-HELPLINE-CMN-RES
$ws.Cells.Item(5, 4).Value = 404
$ws.Cells.Item(5, 5).Value = 406
$ws.Cells.Item(5, 9).Value = 8809
$ws.Cells.Item(5, 10).Value = 12

# Row 6: COV-HELPLINE-FO-CT
$ws.Cells.Item(6, 2).Value = 268
$ws.Cells.Item(6, 4).Value = 15
$ws.Cells.Item(6, 5).Value = 15
$ws.Cells.Item(6, 9).Value = 6000

# Row 7: COV-INFOLINE-PRE-RES
$ws.Cells.Item(7, 2).Value = 264
$ws.Cells.Item(7, 4).Value = 9
$ws.Cells.Item(7, 5).Value = 10
$ws.Cells.Item(7, 8).Value = 10
$ws.Cells.Item(7, 9).Value = 2222

# Row 8: COV-INFOPROV-MOBILE-CMN-RES
$ws.Cells.Item(8, 2).Value = 24
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(8, 4).Value = 498
$ws.Cells.Item(8, 5).Value = 507
$ws.Cells.Item(8, 6).Value = 4
$ws.Cells.Item(8, 8).Value = 79
$ws.Cells.Item(8, 9).Value = 8786

# Row 9: COV-INFOPROVISIONING-CMN-RES
$ws.Cells.Item(9, 2).Value = 775
$ws.Cells.Item(9, 4).Value = 76
$ws.Cells.Item(9, 5).Value = 76
$ws.Cells.Item(9, 9).Value = 8971

# Row 10: COV-INFOPROVISIONING-RES
$ws.Cells.Item(10, 2).Value = 323
$ws.Cells.Item(10, 4).Value = 101
$ws.Cells.Item(10, 5).Value = 116
$ws.Cells.Item(10, 6).Value = 12
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 2722
$ws.Cells.Item(10, 9).Value = 16981
$ws.Cells.Item(10, 10).Value = 5

# Row 11: COV-MIGRAZIONE-35
$ws.Cells.Item(11, 2).Value = 1012
$ws.Cells.Item(11, 4).Value = 465
$ws.Cells.Item(11, 5).Value = 487
$ws.Cells.Item(11, 6).Value = 22
$ws.Cells.Item(11, 7).Value = 3
$ws.Cells.Item(11, 8).Value = 4245
$ws.Cells.Item(11, 9).Value = 42666
$ws.Cells.Item(11, 10).Value = 21

# Row 12: COV-MIGRAZIONE-37
$ws.Cells.Item(12, 2).Value = 842
$ws.Cells.Item(12, 4).Value = 75
$ws.Cells.Item(12, 5).Value = 78
$ws.Cells.Item(12, 6).Value = 5
$ws.Cells.Item(12, 8).Value = 342
$ws.Cells.Item(12, 9).Value = 25873
$ws.Cells.Item(12, 10).Value = 1

# Row 13: COV-MOBILE MNP
$ws.Cells.Item(13, 2).Value = 264
$ws.Cells.Item(13, 3).Value = 5
$ws.Cells.Item(13, 4).Value = 432
$ws.Cells.Item(13, 5).Value = 469
$ws.Cells.Item(13, 6).Value = 34
$ws.Cells.Item(13, 7).Value = 3
$ws.Cells.Item(13, 8).Value = 2541
$ws.Cells.Item(13, 9).Value = 5752
$ws.Cells.Item(13, 10).Value = 11

# Row 14: COV-MOBILE-CMN-RES
$ws.Cells.Item(14, 4).Value = 369
$ws.Cells.Item(14, 5).Value = 371
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 1
$ws.Cells.Item(14, 8).Value = 27
$ws.Cells.Item(14, 9).Value = 9245
$ws.Cells.Item(14, 10).Value = 7

# Row 15: COV-MOBILE-FO-CT
$ws.Cells.Item(15, 2).Value = 266

# Row 16: COV-MSK-GESTIONE-COMM-RES
$ws.Cells.Item(16, 2).Value = 85
$ws.Cells.Item(16, 3).Value = 13
$ws.Cells.Item(16, 4).Value = 710
$ws.Cells.Item(16, 5).Value = 1007
$ws.Cells.Item(16, 6).Value = 56
$ws.Cells.Item(16, 7).Value = 9
$ws.Cells.Item(16, 8).Value = 3121
$ws.Cells.Item(16, 9).Value = 10049
$ws.Cells.Item(16, 10).Value = 6
$ws.Cells.Item(16, 11).Value = 219
$ws.Cells.Item(16, 12).Value = 2607

# Row 17: COV-MSK-HELPLINE-RES
$ws.Cells.Item(17, 2).Value = 66
$ws.Cells.Item(17, 4).Value = 21
$ws.Cells.Item(17, 5).Value = 24
$ws.Cells.Item(17, 6).Value = 2
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 833
$ws.Cells.Item(17, 9).Value = 2381
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 1
$ws.Cells.Item(17, 12).Value = 417

# Row 18: COV-MSK-MOBILE-RES
$ws.Cells.Item(18, 2).Value = 1
$ws.Cells.Item(18, 3).Value = 10
$ws.Cells.Item(18, 4).Value = 672
$ws.Cells.Item(18, 5).Value = 815
$ws.Cells.Item(18, 6).Value = 57
$ws.Cells.Item(18, 7).Value = 1
$ws.Cells.Item(18, 8).Value = 700
$ws.Cells.Item(18, 9).Value = 4747
$ws.Cells.Item(18, 10).Value = 4
$ws.Cells.Item(18, 11).Value = 75
$ws.Cells.Item(18, 12).Value = 920

# Row 19: COV-MSK-RES
$ws.Cells.Item(19, 2).Value = 1515
$ws.Cells.Item(19, 4).Value = 251
$ws.Cells.Item(19, 5).Value = 256
$ws.Cells.Item(19, 6).Value = 5
$ws.Cells.Item(19, 8).Value = 597
$ws.Cells.Item(19, 9).Value = 9214

# Row 20: COV-NOTTE
$ws.Cells.Item(20, 2).Value = 1393

# Row 21: COV-TESTING-RES
$ws.Cells.Item(21, 2).Value = 489
$ws.Cells.Item(21, 4).Value = 93
$ws.Cells.Item(21, 5).Value = 94
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 1
$ws.Cells.Item(21, 8).Value = 1250
$ws.Cells.Item(21, 9).Value = 30003
$ws.Cells.Item(21, 10).Value = 2

# Row 22: COV-VIP-RES
$ws.Cells.Item(22, 2).Value = 191

# Row 23: COV-WELCOME-CMN-RES
$ws.Cells.Item(23, 4).Value = 13
$ws.Cells.Item(23, 5).Value = 13
$ws.Cells.Item(23, 9).Value = 5385

# Row 24: COV-WELCOME-RES
$ws.Cells.Item(24, 2).Value = 999
$ws.Cells.Item(24, 4).Value = 37
$ws.Cells.Item(24, 5).Value = 39
$ws.Cells.Item(24, 6).Value = 2
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = 769
$ws.Cells.Item(24, 9).Value = 13590
$ws.Cells.Item(24, 10).Value = 1

# Row 25: no forecast Condo
$ws.Cells.Item(25, 2).Value = 78
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = 53
$ws.Cells.Item(25, 5).Value = 104
$ws.Cells.Item(25, 6).Value = 8
$ws.Cells.Item(25, 8).Value = 769
$ws.Cells.Item(25, 9).Value = 1887
$ws.Cells.Item(25, 11).Value = 43
$ws.Cells.Item(25, 12).Value = 4135

# Row 26: no forecast Retention
$ws.Cells.Item(26, 2).Value = 137
$ws.Cells.Item(26, 4).Value = 63
$ws.Cells.Item(26, 5).Value = 70
$ws.Cells.Item(26, 6).Value = 5
$ws.Cells.Item(26, 8).Value = 714
$ws.Cells.Item(26, 9).Value = 6349

# Row 27: no forecast vendite
$ws.Cells.Item(27, 2).Value = 15
$ws.Cells.Item(27, 3).Value = 1
$ws.Cells.Item(27, 4).Value = 164
$ws.Cells.Item(27, 5).Value = 201
$ws.Cells.Item(27, 6).Value = 37
$ws.Cells.Item(27, 8).Value = 1841
$ws.Cells.Item(27, 9).Value = 1411
$ws.Cells.Item(27, 10).Value = 3
